# Updates cryptocurrency price (column D) and 1-hour volume change (column E)
# figures on Sheet1, matching the "Updated cryptos list" GitHub Actions commit.
# Values that look like plain decimal numbers (e.g. "0.999", "75.19") are
# entered with a leading apostrophe so Excel stores them as literal text
# (preserving trailing zeros / exact formatting) instead of coercing them
# into numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.661.79"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").Value = "2.421.15"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'317.01"
$ws.Range("E5").Value = "  +4.53%  "
$ws.Range("D6").Value = "'101.23"
$ws.Range("E6").Value = "  +5.94%  "
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +9.65%  "
$ws.Range("D10").Value = "'35.27"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").Value = "'0.0799"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "'18.61"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "2.797.25"
$ws.Range("D16").Value = "2.414.53"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "44.513.64"
$ws.Range("E18").Value = "  +3.31%  "
$ws.Range("D19").Value = "'12.20"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("E21").Value = "  +3.31%  "
$ws.Range("D22").Value = "'68.55"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").Value = "'241.83"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").Value = "'2.27"
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'25.22"
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("D30").Value = "'33.38"
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "'0.125"
$ws.Range("E32").Value = "  +16.43%  "
$ws.Range("D33").Value = "'19.69"
$ws.Range("E33").Value = "  +11.70%  "
$ws.Range("D34").Value = "'5.16"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").Value = "'0.0764"
$ws.Range("E36").Value = "  +5.00%  "
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").Value = "'4.45"
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("D39").Value = "'125.77"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").Value = "'2.85"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  -3.99%  "
$ws.Range("D43").Value = "'20.92"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D45").Value = "1.939.64"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("E47").Value = "  +7.32%  "
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("E49").Value = "  +15.56%  "
$ws.Range("D50").Value = "'75.19"
$ws.Range("E50").Value = "  +5.32%  "
$ws.Range("D51").Value = "'53.59"
$ws.Range("E51").Value = "  +4.74%  "
